# Auto-committed on 2023/09/15 週五 17:07:32.90
# CoreAcMain.xlsx - DBD sheet: rename a field's data type from DATE to TIMESTAMP
# (cells D20 and D22, "CreateDate"/"LastUpdate" rows), and update the
# selection left behind in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Field type for CreateDate (row 20) and LastUpdate (row 22) changes from DATE to TIMESTAMP
$ws.Range("D20").Value = "TIMESTAMP"
$ws.Range("D22").Value = "TIMESTAMP"

# Leave the sheet active with the same selection/scroll state recorded in the saved file
$ws.Activate() | Out-Null
$ws.Range("A4").Select() | Out-Null
$ws.Range("D22").Select() | Out-Null
